$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E5").Value = 8
$ws.Range("E7").Value = 16
$ws.Range("E15").Value = 176
$ws.Range("E17").Value = 140
$ws.Range("F17").Value = 73
$ws.Range("H17").Value = 105
$ws.Range("E18").Value = 133
$ws.Range("F18").Value = 62
$ws.Range("H18").Value = 99
$ws.Range("E27").Value = 16
$ws.Range("F27").Value = 10
$ws.Range("H27").Value = 14
$ws.Range("F33").Value = 15
$ws.Range("H33").Value = 27
$ws.Range("E35").Value = 14
$ws.Range("E36").Value = 120
$ws.Range("E39").Value = 27
$ws.Range("E40").Value = 28
$ws.Range("F40").Value = 19
$ws.Range("H40").Value = 22
$ws.Range("E41").Value = 51
$ws.Range("E48").Value = 41
$ws.Range("E60").Value = 22
$ws.Range("E72").Value = 48
$ws.Range("F72").Value = 29
$ws.Range("H72").Value = 40
$ws.Range("E76").Value = 58
$ws.Range("F76").Value = 23
$ws.Range("H76").Value = 40
$ws.Range("E79").Value = 43
$ws.Range("F79").Value = 21
$ws.Range("H79").Value = 32
$ws.Range("E82").Value = 18
$ws.Range("F82").Value = 7
$ws.Range("H82").Value = 13
$ws.Range("E88").Value = 33
$ws.Range("F88").Value = 20
$ws.Range("H88").Value = 28
